$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Nome: / " Computação..." / same) is never touched by this edit
# and already carries the canonical per-column look (font/alignment), so
# we use it as the formatting donor for every cell we (re)write below.
# This sidesteps two quirks of brand-new cells:
#   1) a string that looks like a date (e.g. "01/01/2020") would otherwise
#      get reinterpreted as a date serial number;
#   2) column B's <col> style definition overlaps with column A's, so a
#      freshly-created B-column cell can otherwise inherit column A's
#      style instead of its own.
function Set-CellText($range, $value, $styleSource) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $styleSource.Copy()
    $range.PasteSpecial(-4122)
    $excel.CutCopyMode = 0
}

$colA = $ws.Range("A3")
$colB = $ws.Range("B3")
$colC = $ws.Range("C3")

# Row 10: Objetivos -> now holds the "519033 - Carlos Yujiro Shigue" text
Set-CellText $ws.Range("A10") "Objetivos:" $colA
Set-CellText $ws.Range("B10") "519033 - Carlos Yujiro Shigue" $colB
Set-CellText $ws.Range("C10") "519033 - Carlos Yujiro Shigue" $colC
$ws.Rows.Item(10).RowHeight = 60

# Row 11: Objectives (text unchanged, just make sure everything is intact)
Set-CellText $ws.Range("A11") "Objectives:" $colA
$ws.Rows.Item(11).RowHeight = 60

# Row 12: Docentes responsáveis (text unchanged)
Set-CellText $ws.Range("A12") "Docentes responsáveis:" $colA
$ws.Rows.Item(12).AutoFit()

# Row 13: now "Programa resumido:" / "01/01/2020"
Set-CellText $ws.Range("A13") "Programa resumido:" $colA
Set-CellText $ws.Range("B13") "01/01/2020" $colB
Set-CellText $ws.Range("C13") "01/01/2020" $colC
$ws.Rows.Item(13).RowHeight = 60

# Row 14: now only "Short syllabus:" in column A
$ws.Range("B14:C14").Clear()
Set-CellText $ws.Range("A14") "Short syllabus:" $colA
$ws.Rows.Item(14).RowHeight = 60

# Row 15: now "Programa:" / "519033 - Carlos Yujiro Shigue"
Set-CellText $ws.Range("A15") "Programa:" $colA
Set-CellText $ws.Range("B15") "519033 - Carlos Yujiro Shigue" $colB
Set-CellText $ws.Range("C15") "519033 - Carlos Yujiro Shigue" $colC
$ws.Rows.Item(15).RowHeight = 120

# Row 16: now "Syllabus:"
Set-CellText $ws.Range("A16") "Syllabus:" $colA
$ws.Rows.Item(16).RowHeight = 120

# Row 17: now only "Avaliação:" in column A
$ws.Range("B17:C17").Clear()
Set-CellText $ws.Range("A17") "Avaliação:" $colA
$ws.Rows.Item(17).AutoFit()

# Row 18: now "Método:" / "1176388 - Luiz Tadeu Fernandes Eleno"
Set-CellText $ws.Range("A18") "Método:" $colA
Set-CellText $ws.Range("B18") "1176388 - Luiz Tadeu Fernandes Eleno" $colB
Set-CellText $ws.Range("C18") "1176388 - Luiz Tadeu Fernandes Eleno" $colC
$ws.Rows.Item(18).RowHeight = 60

# Row 19: now "Critério:" / the "Aulas expositivas..." text
Set-CellText $ws.Range("A19") "Critério:" $colA
Set-CellText $ws.Range("B19") "Aulas expositivas e em laboratório computacional, trabalhos e exercícios comentados." $colB
Set-CellText $ws.Range("C19") "Aulas expositivas e em laboratório computacional, trabalhos e exercícios comentados." $colC
$ws.Rows.Item(19).RowHeight = 60

# Row 20: now "Norma de recuperação:" / "Média aritmética..." text
Set-CellText $ws.Range("A20") "Norma de recuperação:" $colA
Set-CellText $ws.Range("B20") "Média aritmética de trabalhos propostos ao longo do curso." $colB
Set-CellText $ws.Range("C20") "Média aritmética de trabalhos propostos ao longo do curso." $colC
$ws.Rows.Item(20).RowHeight = 60

# Row 21: now "Bibliografia:" / "Não haverá exame de recuperação."
Set-CellText $ws.Range("A21") "Bibliografia:" $colA
Set-CellText $ws.Range("B21") "Não haverá exame de recuperação." $colB
Set-CellText $ws.Range("C21") "Não haverá exame de recuperação." $colC
$ws.Rows.Item(21).RowHeight = 120

# Rows 22 and 23 no longer exist in the trimmed sheet - delete from the
# bottom up so the row numbers above stay stable while deleting.
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(22).Delete()

Write-Output "done"
